$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 13 (docente name, with no label in column A) is removed; everything
# below shifts up by one row and most rows then already line up with the
# desired labels/heights. Only the B/C text of a handful of rows still needs
# to be corrected afterwards.
$ws.Rows.Item(13).Delete()

$ws.Range("B10").Value = "5840917 - Fabrício Maciel Gomes"
$ws.Range("C10").Value = "5840917 - Fabrício Maciel Gomes"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2021"
$ws.Range("C15").Value = "01/01/2021"

$ws.Range("B18").Value = "5840917 - Fabrício Maciel Gomes"
$ws.Range("C18").Value = "5840917 - Fabrício Maciel Gomes"

$ws.Range("B19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("C19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

$ws.Range("B20").Value = "NF≥ 5,0."
$ws.Range("C20").Value = "NF≥ 5,0."

$ws.Range("B21").Value = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação."
$ws.Range("C21").Value = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação."
